# Update cryptos list: price (D) and 1h volume/change (E) columns for rows 2-50,
# and rename row 51 from "dogwifhat" to "Polygon" with new link/price/change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $text) {
    # Force the cell to remain a text value (matches source data where
    # numeric-looking price strings like "1.00" or "0.387" are stored as text)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "57.544.46"
Set-TextValue $ws.Range("E2") "  +1.18%  "
Set-TextValue $ws.Range("D3") "2.363.96"
Set-TextValue $ws.Range("E3") "  +1.00%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "521.13"
Set-TextValue $ws.Range("E5") "  +0.65%  "
Set-TextValue $ws.Range("D6") "136.12"
Set-TextValue $ws.Range("E6") "  +1.51%  "
Set-TextValue $ws.Range("D7") "0.998"
Set-TextValue $ws.Range("E7") "  -0.11%  "
Set-TextValue $ws.Range("D8") "0.540"
Set-TextValue $ws.Range("E8") "  +0.68%  "
Set-TextValue $ws.Range("E9") "  -0.46%  "
Set-TextValue $ws.Range("D10") "5.43"
Set-TextValue $ws.Range("E10") "  +5.17%  "
Set-TextValue $ws.Range("E11") "  -0.71%  "
Set-TextValue $ws.Range("E12") "  -0.05%  "
Set-TextValue $ws.Range("D13") "24.34"
Set-TextValue $ws.Range("E13") "  +1.41%  "
Set-TextValue $ws.Range("D14") "2.789.46"
Set-TextValue $ws.Range("E14") "  +2.07%  "
Set-TextValue $ws.Range("D15") "57.533.74"
Set-TextValue $ws.Range("E15") "  +1.35%  "
Set-TextValue $ws.Range("E16") "  +0.46%  "
Set-TextValue $ws.Range("D17") "2.387.78"
Set-TextValue $ws.Range("E17") "  +2.10%  "
Set-TextValue $ws.Range("E18") "  +0.92%  "
Set-TextValue $ws.Range("D19") "330.92"
Set-TextValue $ws.Range("E19") "  +3.00%  "
Set-TextValue $ws.Range("D20") "4.24"
Set-TextValue $ws.Range("E20") "  -1.08%  "
Set-TextValue $ws.Range("E21") "  +1.15%  "
Set-TextValue $ws.Range("E22") "  +0.24%  "
Set-TextValue $ws.Range("D23") "61.32"
Set-TextValue $ws.Range("E23") "  +0.14%  "
Set-TextValue $ws.Range("D24") "8.80"
Set-TextValue $ws.Range("E24") "  +14.66%  "
Set-TextValue $ws.Range("E25") "  +4.59%  "
Set-TextValue $ws.Range("D26") "0.996"
Set-TextValue $ws.Range("E26") "  +0.33%  "
Set-TextValue $ws.Range("E27") "  +12.04%  "
Set-TextValue $ws.Range("D28") "0.0₃0749"
Set-TextValue $ws.Range("E28") "  +1.72%  "
Set-TextValue $ws.Range("D29") "169.58"
Set-TextValue $ws.Range("E29") "  -1.26%  "
Set-TextValue $ws.Range("E30") "  +1.75%  "
Set-TextValue $ws.Range("D32") "18.61"
Set-TextValue $ws.Range("E32") "  +1.34%  "
Set-TextValue $ws.Range("E33") "  +0.04%  "
Set-TextValue $ws.Range("D34") "1.32"
Set-TextValue $ws.Range("E34") "  +3.88%  "
Set-TextValue $ws.Range("D35") "0.996"
Set-TextValue $ws.Range("E35") "  -0.34%  "
Set-TextValue $ws.Range("D36") "0.924"
Set-TextValue $ws.Range("E36") "  -2.07%  "
Set-TextValue $ws.Range("D37") "4.07"
Set-TextValue $ws.Range("E37") "  +0.98%  "
Set-TextValue $ws.Range("E38") "  +7.79%  "
Set-TextValue $ws.Range("D39") "38.68"
Set-TextValue $ws.Range("E39") "  +2.97%  "
Set-TextValue $ws.Range("D40") "151.07"
Set-TextValue $ws.Range("E40") "  +7.25%  "
Set-TextValue $ws.Range("D41") "0.389"
Set-TextValue $ws.Range("E41") "  +1.86%  "
Set-TextValue $ws.Range("D42") "3.68"
Set-TextValue $ws.Range("E42") "  +2.06%  "
Set-TextValue $ws.Range("E43") "  +3.19%  "
Set-TextValue $ws.Range("D44") "284.13"
Set-TextValue $ws.Range("E44") "  +2.87%  "
Set-TextValue $ws.Range("D45") "0.0944"
Set-TextValue $ws.Range("E45") "  +1.55%  "
Set-TextValue $ws.Range("E46") "  +0.10%  "
Set-TextValue $ws.Range("D47") "0.565"
Set-TextValue $ws.Range("E47") "  +0.61%  "
Set-TextValue $ws.Range("D48") "18.35"
Set-TextValue $ws.Range("E48") "  +6.29%  "
Set-TextValue $ws.Range("E49") "  +2.26%  "
Set-TextValue $ws.Range("D50") "17.89"
Set-TextValue $ws.Range("E50") "  +5.67%  "

# Row 51: coin changed from dogwifhat to Polygon
Set-TextValue $ws.Range("B51") "Polygon"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D51") "0.387"
Set-TextValue $ws.Range("E51") "  +1.52%  "

